$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.42
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 8.5
$ws.Range("J2").Value = 2
$ws.Range("L2").Value = 8
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("AA2").Value = 2.38
$ws.Range("AB2").Value = 1.53
$ws.Range("AE2").Value = 9
$ws.Range("AF2").Value = 9
$ws.Range("AI2").Value = 8
$ws.Range("AM2").Value = 15
$ws.Range("AR2").Value = 67
$ws.Range("G3").Value = 3.25
$ws.Range("I3").Value = 2.25
$ws.Range("J3").Value = 4
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 7.5
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 2.35
$ws.Range("T3").Value = 1.57
$ws.Range("U3").Value = 3.15
$ws.Range("V3").Value = 1.34
$ws.Range("AA3").Value = 2
$ws.Range("AB3").Value = 1.73
$ws.Range("AE3").Value = 12
$ws.Range("AF3").Value = 34
$ws.Range("AG3").Value = 29
$ws.Range("AI3").Value = 7.5
$ws.Range("AK3").Value = 17
$ws.Range("AM3").Value = 6.5
$ws.Range("AN3").Value = 10
$ws.Range("AR3").Value = 34
$ws.Range("AS3").Value = 1000
